$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.301.74'
$ws.Range("E2").Value = '  -2.17%  '

$ws.Range("D3").Value = '3.709.95'
$ws.Range("E3").Value = '  -2.99%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.39%  '

$ws.Range("D7").Value = '3.707.66'
$ws.Range("E7").Value = '  -3.09%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.25%  '

$ws.Range("E12").Value = '  -3.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.04%  '

$ws.Range("E14").Value = '  -2.38%  '

$ws.Range("D15").Value = '4.338.72'
$ws.Range("E15").Value = '  -3.18%  '

$ws.Range("D16").Value = '3.718.10'
$ws.Range("E16").Value = '  -3.26%  '

$ws.Range("D17").Value = '67.383.94'
$ws.Range("E17").Value = '  -2.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.79%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '487.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.727'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000143'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.02'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.15%  '

$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("E30").Value = '  -1.54%  '

$ws.Range("E31").Value = '  -5.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.61%  '

$ws.Range("D34").Value = '3.859.20'
$ws.Range("E34").Value = '  -2.93%  '

$ws.Range("E35").Value = '  -3.61%  '

$ws.Range("D36").Value = '3.660.69'
$ws.Range("E36").Value = '  -2.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.997'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.16%  '

$ws.Range("E39").Value = '  -2.70%  '

$ws.Range("E40").Value = '  -4.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.321'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '48.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '427.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.15%  '

$ws.Range("E44").Value = '  -1.39%  '

$ws.Range("E45").Value = '  -4.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.31%  '

$ws.Range("D50").Value = '2.754.18'
$ws.Range("E50").Value = '  -4.40%  '

$ws.Range("E51").Value = '  -2.08%  '
